$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL) swaps: safe to set directly ---
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

# --- Numeric-look-alike text columns (Price / Volume%) ---
# These values (e.g. "1.001", "30.069.76") would be auto-coerced to
# numbers by plain .Value assignment since they parse as numeric literals.
# Force text storage by switching the cell to the "@" (Text) number format
# before assigning, then restore the original (General/default) formatting
# by pasting formats from an untouched, same-style text cell (column B),
# so the cell style index is not left modified.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.069.76'
$ws.Range("E2").Value = '  +5.42%  '
$ws.Range("D3").Value = '1.923.05'
$ws.Range("E3").Value = '  +2.62%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.73%  '
$ws.Range("D5").Value = '325.68'
$ws.Range("E5").Value = '  +3.22%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("D7").Value = '0.5159'
$ws.Range("E7").Value = '  +1.64%  '
$ws.Range("D8").Value = '0.3998'
$ws.Range("E8").Value = '  +2.62%  '
$ws.Range("D9").Value = '0.08468'
$ws.Range("D10").Value = '42.93'
$ws.Range("E10").Value = '  +2.73%  '
$ws.Range("D11").Value = '1.122'
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").Value = '6.333'
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").Value = '21.17'
$ws.Range("E13").Value = '  +3.87%  '
$ws.Range("D14").Value = '1.919.29'
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").Value = '7.338'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("D17").Value = '94.26'
$ws.Range("E17").Value = '  +3.39%  '
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").Value = '17.99'
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("D22").Value = '6.059'
$ws.Range("E22").Value = '  +2.16%  '
$ws.Range("D23").Value = '30.072.66'
$ws.Range("E23").Value = '  +5.30%  '
$ws.Range("D24").Value = '11.18'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("D26").Value = '2.141.79'
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("D27").Value = '160.11'
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("D29").Value = '2.464'
$ws.Range("E29").Value = '  +4.61%  '
$ws.Range("D30").Value = '129.02'
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("D31").Value = '1.079'
$ws.Range("E31").Value = '  +3.52%  '
$ws.Range("D32").Value = '0.1058'
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").Value = '6.073'
$ws.Range("E33").Value = '  +4.85%  '
$ws.Range("D34").Value = '3.660'
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("D35").Value = '0.02500'
$ws.Range("E35").Value = '  +2.00%  '
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").Value = '0.2221'
$ws.Range("E37").Value = '  +2.81%  '
$ws.Range("D38").Value = '1.244'
$ws.Range("E38").Value = '  +4.38%  '
$ws.Range("D39").Value = '5.197'
$ws.Range("E39").Value = '  +2.73%  '
$ws.Range("D40").Value = '8.988'
$ws.Range("E40").Value = '  +1.26%  '
$ws.Range("D41").Value = '0.6532'
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D43").Value = '11.41'
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("D44").Value = '0.6132'
$ws.Range("E44").Value = '  +1.60%  '
$ws.Range("D45").Value = '13.14'
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("D46").Value = '3.741'
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("D47").Value = '2.057'
$ws.Range("E47").Value = '  +2.27%  '
$ws.Range("D48").Value = '1.244'
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("D49").Value = '125.44'
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("D50").Value = '79.34'
$ws.Range("E50").Value = '  +3.51%  '
$ws.Range("E51").Value = '  -2.57%  '

# Restore default (General) style/format on the whole D:E block in one shot
# by copying formats from the never-touched, default-styled B column.
$ws.Range("B2:B51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4122)
$ws.Range("B2:B51").Copy()
$ws.Range("E2:E51").PasteSpecial(-4122)
